{"js": "// Apply review-comment edits to the DME 2.34.0 release notes document.\n\n// ---------------------------------------------------------------------\n// 1) HPCDATAMGM-1780, 1800 paragraph: \"screen\" -> \"page\", add \"other \"\n//    and \"own \", and add a comma after \"For details\".\n// ---------------------------------------------------------------------\n{\n  const r1 = context.document.body.search(\n    \"screen of the DME web application\",\n    { matchCase: true, matchWholeWord: false }\n  );\n  r1.load(\"items\");\n  await context.sync();\n  if (r1.items.length > 0) {\n    r1.items[0].insertText(\"page of the DME web application\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n{\n  const r2 = context.document.body.search(\n    \"notification subscriptions for users in their DOC\",\n    { matchCase: true }\n  );\n  r2.load(\"items\");\n  await context.sync();\n  if (r2.items.length > 0) {\n    r2.items[0].insertText(\n      \"notification subscriptions for other users in their DOC\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n}\n\n{\n  const r3 = context.document.body.search(\n    \"could only add/update their subscriptions via the GUI\",\n    { matchCase: true }\n  );\n  r3.load(\"items\");\n  await context.sync();\n  if (r3.items.length > 0) {\n    r3.items[0].insertText(\n      \"could only add/update their own subscriptions via the GUI\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n}\n\n{\n  const r4 = context.document.body.search(\n    \"irrespective of their role. For details refer to\",\n    { matchCase: true }\n  );\n  r4.load(\"items\");\n  await context.sync();\n  if (r4.items.length > 0) {\n    r4.items[0].insertText(\n      \"irrespective of their role. For details, refer to\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n}\n\n// ---------------------------------------------------------------------\n// 2) HPCDATAMGM-1795 paragraph: \"refer section 5.31\" -> \"refer to section 5.31\"\n// ---------------------------------------------------------------------\n{\n  const r5 = context.document.body.search(\n    \"For details, refer section 5.31\",\n    { matchCase: true }\n  );\n  r5.load(\"items\");\n  await context.sync();\n  if (r5.items.length > 0) {\n    r5.items[0].insertText(\n      \"For details, refer to section 5.31\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n}\n\n// ---------------------------------------------------------------------\n// 3) HPCDATAMGM-1791 paragraph: drop \"single row or the \" and italicize\n//    \"Select All\".\n// ---------------------------------------------------------------------\n{\n  const r6 = context.document.body.search(\n    \"displayed when the single row or the Select All checkbox\",\n    { matchCase: true }\n  );\n  r6.load(\"items\");\n  await context.sync();\n  if (r6.items.length > 0) {\n    r6.items[0].insertText(\n      \"displayed when the Select All checkbox\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n}\n\n{\n  const r7 = context.document.body.search(\"Select All\", { matchCase: true });\n  r7.load(\"items\");\n  await context.sync();\n  if (r7.items.length > 0) {\n    r7.items[0].font.italic = true;\n    await context.sync();\n  }\n}\n\n// ---------------------------------------------------------------------\n// 4) Important Notes paragraph: \"Release.\" -> \"release.\" (lowercase the R)\n//    Keep this a minimal, single-character replacement so surrounding\n//    (non-breaking) whitespace is left untouched.\n// ---------------------------------------------------------------------\n{\n  const r8 = context.document.body.search(\"Release.\", { matchCase: true });\n  r8.load(\"items\");\n  await context.sync();\n  if (r8.items.length > 0) {\n    r8.items[0].insertText(\"release.\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n", "ps1": "# Apply review-comment edits to the DME 2.34.0 release notes document.\n\n$d = $word.ActiveDocument\n$wdFindContinue = 1\n$wdReplaceOne   = 1\n\n# Helper: run a single Find/Replace pass (replaces only the first match,\n# same semantics as a surgical, single-occurrence replacement) over the\n# whole document body. Keeping each replacement minimal (instead of\n# rewriting whole sentences) avoids touching unrelated whitespace, such\n# as the non-breaking spaces already present elsewhere in this document.\nfunction Replace-FirstMatch {\n    param(\n        [string]$FindText,\n        [string]$ReplaceWith,\n        [bool]$MatchCase = $true\n    )\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    return $rng.Find.Execute($FindText, $MatchCase, $true, $false, $false, $false, $true, $wdFindContinue, $false, $ReplaceWith, $wdReplaceOne)\n}\n\n# ---------------------------------------------------------------------\n# 1) HPCDATAMGM-1780, 1800 paragraph: \"screen\" -> \"page\", add \"other \"\n#    and \"own \", and add a comma after \"For details\".\n# ---------------------------------------------------------------------\nReplace-FirstMatch \"screen of the DME web application\" \"page of the DME web application\" | Out-Null\nReplace-FirstMatch \"notification subscriptions for users in their DOC\" \"notification subscriptions for other users in their DOC\" | Out-Null\nReplace-FirstMatch \"could only add/update their subscriptions via the GUI\" \"could only add/update their own subscriptions via the GUI\" | Out-Null\nReplace-FirstMatch \"irrespective of their role. For details refer to\" \"irrespective of their role. For details, refer to\" | Out-Null\n\n# ---------------------------------------------------------------------\n# 2) HPCDATAMGM-1795 paragraph: \"refer section 5.31\" -> \"refer to section 5.31\"\n# ---------------------------------------------------------------------\nReplace-FirstMatch \"For details, refer section 5.31\" \"For details, refer to section 5.31\" | Out-Null\n\n# ---------------------------------------------------------------------\n# 3) HPCDATAMGM-1791 paragraph: drop \"single row or the \" and italicize\n#    \"Select All\".\n# ---------------------------------------------------------------------\nReplace-FirstMatch \"displayed when the single row or the Select All checkbox\" \"displayed when the Select All checkbox\" | Out-Null\n\n$selectAllRange = $d.Content\n$selectAllRange.Find.ClearFormatting()\n$selectAllRange.Find.Text = \"Select All\"\n$selectAllRange.Find.MatchCase = $true\n$selectAllRange.Find.MatchWholeWord = $false\nif ($selectAllRange.Find.Execute()) {\n    $selectAllRange.Font.Italic = 1\n}\n\n# ---------------------------------------------------------------------\n# 4) Important Notes paragraph: \"Release.\" -> \"release.\" (lowercase the R)\n#    Keep this a minimal, single-character replacement so surrounding\n#    (non-breaking) whitespace is left untouched.\n# ---------------------------------------------------------------------\nReplace-FirstMatch \"Release.\" \"release.\" | Out-Null\n"}
